# Adding DataBricks SQL Endpoint cases
# Update the "suites" worksheet:
#   B2 -> new generated frontEnd result file name
#   B3 -> "null" (authentication test result became null)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("suites")

$ws.Range("B2").Value = "/target/frontEndqdALMMFIHLzQrwWs.html"
$ws.Range("B3").Value = "null"
